$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows 251 and 252 hold the most recent weekly price report for
# "Vega Monumental Concepción - Cilantro". A new week of data has come in,
# so the existing rows move down to become historical rows 253/254
# (preserving their current values), and rows 251/252 are updated in
# place with the new date and prices.

# 1) Copy the current (soon-to-be-historical) rows 251:252 down to 253:254.
$ws.Range("A251:R252").Copy()
$ws.Range("A253").PasteSpecial()

# 2) Update row 251 with the new report's date and prices.
$ws.Range("D251").Value = 44911
$ws.Range("K251").Value = 700
$ws.Range("L251").Value = 800
$ws.Range("M251").Value = 750
$ws.Range("P251").Value = 750

# 3) Update row 252 with the new report's date and prices.
$ws.Range("D252").Value = 44911
$ws.Range("K252").Value = 600
$ws.Range("L252").Value = 600
$ws.Range("M252").Value = 600
$ws.Range("P252").Value = 600
